$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.830.78'
$ws.Range("E2").Value = '  +1.05%  '
$ws.Range("D3").Value = '2.602.43'
$ws.Range("E3").Value = '  +0.65%  '
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '557.31'
$ws.Range("E5").Value = '  -1.80%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.78'
$ws.Range("E6").Value = '  -1.31%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.600'
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").Value = '2.626.07'
$ws.Range("E9").Value = '  +1.30%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.67'
$ws.Range("E10").Value = '  -0.03%  '
$ws.Range("E11").Value = '  +1.34%  '
$ws.Range("E12").Value = '  +4.99%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.371'
$ws.Range("E13").Value = '  +8.24%  '
$ws.Range("D14").Value = '3.066.70'
$ws.Range("E14").Value = '  +0.90%  '
$ws.Range("B15").Value = 'WrappedBTC'
$ws.Range("C15").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D15").Value = '59.813.43'
$ws.Range("E15").Value = '  +0.96%  '
$ws.Range("B16").Value = 'Avalanche'
$ws.Range("C16").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '23.34'
$ws.Range("E16").Value = '  +4.12%  '
$ws.Range("E17").Value = '  +0.57%  '
$ws.Range("D18").Value = '2.621.37'
$ws.Range("E18").Value = '  +1.28%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.63'
$ws.Range("E19").Value = '  +2.79%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '343.66'
$ws.Range("E20").Value = '  +2.10%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.73'
$ws.Range("E21").Value = '  +5.43%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.89'
$ws.Range("E22").Value = '  +11.14%  '
$ws.Range("E23").Value = '  +0.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.521'
$ws.Range("E24").Value = '  +15.98%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '62.61'
$ws.Range("E25").Value = '  -2.34%  '
$ws.Range("E26").Value = '  -0.01%  '
$ws.Range("E27").Value = '  -1.28%  '
$ws.Range("E28").Value = '  +5.88%  '
$ws.Range("E29").Value = '  -0.10%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.998'
$ws.Range("E30").Value = '  +0.03%  '
$ws.Range("E31").Value = '  +1.30%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.19'
$ws.Range("E32").Value = '  +2.42%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '159.06'
$ws.Range("E33").Value = '  +0.25%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '19.49'
$ws.Range("E34").Value = '  +2.50%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.17'
$ws.Range("E35").Value = '  +3.28%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.935'
$ws.Range("E36").Value = '  +5.50%  '
$ws.Range("E37").Value = '  +4.84%  '
$ws.Range("B38").Value = 'Stacks'
$ws.Range("C38").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.54'
$ws.Range("E38").Value = '  +2.42%  '
$ws.Range("B39").Value = 'OKB'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '37.77'
$ws.Range("E39").Value = '  +2.15%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.850'
$ws.Range("E40").Value = '  -2.80%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.75'
$ws.Range("E41").Value = '  +2.46%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '298.05'
$ws.Range("E42").Value = '  +1.94%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '141.73'
$ws.Range("E43").Value = '  +14.16%  '
$ws.Range("E44").Value = '  -0.18%  '
$ws.Range("E45").Value = '  +0.75%  '
$ws.Range("E46").Value = '  +0.73%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0241'
$ws.Range("E47").Value = '  +3.70%  '
$ws.Range("E48").Value = '  +0.73%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '10.64'
$ws.Range("E49").Value = '  +0.00%  '
$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '19.24'
$ws.Range("E50").Value = '  +3.66%  '
$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").Value = '2.000.14'
$ws.Range("E51").Value = '  +3.17%  '
